# Add a new record row (row 13) to the "data" sheet, mirroring the
# existing rows 2-12:
#   A13 = ID (uuid)              -> bold/bordered/centered style like col A
#   B13 = created date (text)    -> "2023-01-30"
#   C13 = modified (blank)
#   D13 = name                   -> "taxes2022"
#   E13 = source                 -> "taxes2022.xlsx"
#   F13 = status                 -> "active"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: copy the header-ish record style used by A2:A12 (bold,
#     bordered, centered) onto the new A13 cell before setting its value.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A13").Value = "bccf28a2-848e-4783-ad20-4da0c98f73cd"

# --- Column B: the "created" column stores plain date-like text (as in
#     B6:B12), not a real date value, so force a text number format
#     before assigning, then drop back to the plain/default style used
#     by the rest of the column so no stray formatting is left behind.
$ws.Range("C13").NumberFormat = "@"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "2023-01-30"
$ws.Range("B13").Style = $ws.Range("B2").Style

# --- Column C: "modified" is left blank for this record, same as every
#     other row in the sheet; just align its formatting with the rest
#     of the (unstyled) column.
$ws.Range("C13").Style = $ws.Range("C2").Style

# --- Columns D-F: plain text values.
$ws.Range("D13").Value = "taxes2022"
$ws.Range("E13").Value = "taxes2022.xlsx"
$ws.Range("F13").Value = "active"
